{"js": "// The paragraph (on f.121v) reads, in the transcription markup:\n//   ... ou une <ms>petite <bp>joinctee</bp></ms> puys mesle ...\n// The edit removes the <bp>...</bp> markup around \"joinctee\" while\n// keeping the plain text, merging \"petite \" and \"joinctee\" into a\n// single run, and collapsing \"</bp></ms>\" down to just \"</ms>\".\n//\n// \"petite\" / \"joinctee\" / the \"<bp>\" right before \"joinctee\" are all\n// unique in the document, so we can target them unambiguously with\n// Word's search.\n\nconst body = context.document.body;\n\n// Step 1: merge \"petite \" + \"<bp>\" + \"joinctee\" into a single run of\n// plain text \"petite joinctee\" (keeps the formatting of the first run\n// in the matched range, i.e. the Arial/black/22 run \"petite \" already\n// carried).\nconst openTagMatches = body.search(\"petite <bp>joinctee\", { matchCase: true });\nopenTagMatches.load(\"text\");\nawait context.sync();\n\nif (openTagMatches.items.length !== 1) {\n  throw new Error(\n    \"expected exactly one match for 'petite <bp>joinctee', found \" +\n      openTagMatches.items.length\n  );\n}\nopenTagMatches.items[0].insertText(\"petite joinctee\", \"Replace\");\nawait context.sync();\n\n// Step 2: remove the closing \"</bp>\" that immediately precedes \"</ms>\"\n// right after \"joinctee\", leaving just \"</ms>\" in its own (Courier\n// New / blue) run.\nconst closingScope = body.search(\"joinctee</bp></ms>\", { matchCase: true });\nclosingScope.load(\"text\");\nawait context.sync();\n\nif (closingScope.items.length !== 1) {\n  throw new Error(\n    \"expected exactly one match for 'joinctee</bp></ms>', found \" +\n      closingScope.items.length\n  );\n}\n\nconst closeTagMatches = closingScope.items[0].search(\"</bp>\", { matchCase: true });\ncloseTagMatches.load(\"text\");\nawait context.sync();\n\nif (closeTagMatches.items.length !== 1) {\n  throw new Error(\n    \"expected exactly one '</bp>' inside the matched scope, found \" +\n      closeTagMatches.items.length\n  );\n}\n\ncloseTagMatches.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n", "ps1": "# The paragraph (on f.121v) reads, in the transcription markup:\n#   ... ou une <ms>petite <bp>joinctee</bp></ms> puys mesle ...\n# The edit removes the <bp>...</bp> markup around \"joinctee\" while\n# keeping the plain text, merging \"petite \" and \"joinctee\" into a\n# single run, and collapsing \"</bp></ms>\" down to just \"</ms>\".\n#\n# \"petite\" / \"joinctee\" / the \"<bp>\" right before \"joinctee\" are all\n# unique in the document, so we can target them unambiguously with\n# Word's Find.\n\n$d = $word.ActiveDocument\n\n# Step 1: merge \"petite \" + \"<bp>\" + \"joinctee\" into a single run of\n# plain text \"petite joinctee\" (Range.Text replaces the whole matched\n# span and keeps the formatting carried by its first run, i.e. the\n# Arial/black/22 \"petite \" run).\n$openRange = $d.Content\n$openRange.Find.ClearFormatting()\n$openRange.Find.Text = \"petite <bp>joinctee\"\n$openRange.Find.MatchCase = $true\n$openRange.Find.MatchWildcards = $false\n$found1 = $openRange.Find.Execute()\nif (-not $found1) {\n    throw \"could not find 'petite <bp>joinctee'\"\n}\n$openRange.Text = \"petite joinctee\"\n\n# Step 2: remove the closing \"</bp>\" that immediately precedes \"</ms>\"\n# right after \"joinctee\", leaving just \"</ms>\" in its own (Courier\n# New / blue) run.\n$closeScope = $d.Content\n$closeScope.Find.ClearFormatting()\n$closeScope.Find.Text = \"joinctee</bp></ms>\"\n$closeScope.Find.MatchCase = $true\n$closeScope.Find.MatchWildcards = $false\n$found2 = $closeScope.Find.Execute()\nif (-not $found2) {\n    throw \"could not find 'joinctee</bp></ms>'\"\n}\n\n# Duplicate the narrow scope so we can run an independent Find bounded\n# to it (searching inside the already-collapsed \"found\" range needs a\n# fresh Range object).\n$closeTag = $closeScope.Duplicate\n$closeTag.Find.ClearFormatting()\n$closeTag.Find.Text = \"</bp>\"\n$closeTag.Find.MatchCase = $true\n$closeTag.Find.MatchWildcards = $false\n$found3 = $closeTag.Find.Execute()\nif (-not $found3) {\n    throw \"could not find '</bp>' before '</ms>'\"\n}\n$closeTag.Text = \"\"\n"}
